$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two data rows (record 111609176 on row 2, record 111609174 on row 3)
# traded places in the sheet. Only touch the columns whose values actually
# differ between the two rows; columns holding identical values on both
# rows are left alone (swapping them would be a no-op, and re-writing
# some of them, e.g. the text dates in Y/Z/AA/AB, would risk Excel
# re-interpreting the literal text as a real date/number).

$numericCols = @("A","B","E","Q","R")
$textCols    = @("D","F","G","H")

foreach ($col in $numericCols) {
    $addr2 = $col + "2"
    $addr3 = $col + "3"
    $v2 = $ws.Range($addr2).Value2
    $v3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value = $v3
    $ws.Range($addr3).Value = $v2
}

foreach ($col in $textCols) {
    $addr2 = $col + "2"
    $addr3 = $col + "3"
    $v2 = $ws.Range($addr2).Value2
    $v3 = $ws.Range($addr3).Value2
    $ws.Range($addr2).Value = $v3
    $ws.Range($addr3).Value = $v2
}

# Column I: "30" (text) moves from row 2 to row 3; row 2 becomes blank.
# Force text type via a temporary "@" number format so "30" isn't stored
# as a numeric value, then clear the format again so no stray style is
# left behind on the cell.
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "30"
$ws.Range("I3").ClearFormats()
$ws.Range("I2").Value = ""

# Column AO: habitat note only present on row 3 before the edit; it moves
# to row 2, and row 3 loses it.
$ws.Range("AO2").Value = "mossig silverl" + [char]0x00E5 + "ga av tall"
$ws.Range("AO3").Value = ""
